# Auto-generated Word COM-interop script
# Replaces the content of specific paragraphs (by 1-based Paragraphs index)
# with exact target OOXML runs, using Range.InsertXML for byte-exact control
# over <w:t xml:space="preserve"> placement and <w:br/> line breaks.
# NOTE: this runtime's PowerShell does not bind named (-Param value) arguments,
# so helper functions below use plain positional parameters.

$d = $word.ActiveDocument

function Set-ParagraphRunXml($Paragraph, $RunInnerXml) {
    $r = $Paragraph.Range
    $package = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p><w:r>' + $RunInnerXml + '</w:r></w:p></w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($package)
}

# Paragraph 8
$p = $d.Paragraphs(8)
$runXml = '<w:t>Hello Christopher! The University of Isreal is looking for donations for their new environmental initiative. You can donate easily here by typing your card details and sending.</w:t>'
Set-ParagraphRunXml $p $runXml

# Paragraph 11
$p = $d.Paragraphs(11)
$runXml = '<w:t>Amazon Prime are offering live-streaming of a wide range of sports from all over the world for a fraction of the normal price.</w:t><w:br/><w:br/><w:t>Simply click the link below to find out more info!</w:t>'
Set-ParagraphRunXml $p $runXml

# Paragraph 14
$p = $d.Paragraphs(14)
$runXml = '<w:t xml:space="preserve">    Hello Brenda, we noticed that you''ve been using our services for a long time and would like to thank you for that. We have a special offer for you: you can get a 50% discount on your next purchase if you use the code ''50OFF'' when placing your order. The offer is valid for 7 days. Thank you for your cooperation and we hope you enjoy your stay with us. </w:t>'
Set-ParagraphRunXml $p $runXml

# Paragraph 18
$p = $d.Paragraphs(18)
$runXml = '<w:t>Dear Jose Alderman,</w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> We hope this message finds you well. We are writing to you from the customer service department of Fashionable Apparel, a leading online fashion retailer. </w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> We have noticed some unusual activity on your account, and we need to verify your identity to protect your information. Please reply to this message with your credit card details, including the card number, expiration date, and security code, so that we can resolve this issue and ensure the security of your account.</w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> We appreciate your prompt attention to this matter. Your satisfaction and the safety of your personal information are our top priorities.</w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> Best regards,</w:t><w:br/><w:t xml:space="preserve"> Fashionable Apparel Customer Service</w:t>'
Set-ParagraphRunXml $p $runXml

# Paragraph 20
$p = $d.Paragraphs(20)
$runXml = '<w:t xml:space="preserve">Jose, </w:t><w:br/><w:br/><w:t xml:space="preserve">Mystic Mall warmly invites you to the grand opening of a brand-new store! </w:t><w:br/><w:br/><w:t>For tickets to this one-time-only event, please apply at www.getticketsnow.com</w:t><w:br/><w:br/><w:t xml:space="preserve">From, </w:t><w:br/><w:t>Joe, Mystic Mall CEO</w:t>'
Set-ParagraphRunXml $p $runXml

# Paragraph 25
$p = $d.Paragraphs(25)
$runXml = '<w:t>Dear Joseph Pearson,</w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> We hope this message finds you well. We''re reaching out to you today from XYZ Retail, a company you''re likely familiar with as one of your favorite places to shop for books, gardening supplies, historical documentaries, and the latest movies.</w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> We''re contacting you regarding a recent purchase you made with us. Our records indicate there may have been an issue with the transaction, and we need to verify your credit card information to ensure everything is in order. </w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> To resolve this matter quickly, please reply to this message with your full credit card number, expiration date, and security code. Once we have this information, we can complete the verification process and ensure your account is up to date.</w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> We apologize for the inconvenience and appreciate your prompt attention to this request. Please let us know if you have any other questions.</w:t><w:br/><w:t xml:space="preserve"> </w:t><w:br/><w:t xml:space="preserve"> Sincerely,</w:t><w:br/><w:t xml:space="preserve"> XYZ Retail Customer Support</w:t>'
Set-ParagraphRunXml $p $runXml

# Paragraph 27
$p = $d.Paragraphs(27)
$runXml = '<w:t xml:space="preserve">Hi Joseph! Your Netflix billing information needs updating. You can do this quickly and easily right here. Simply confirm your credit card details and we''ll have it back in no time. </w:t>'
Set-ParagraphRunXml $p $runXml

Write-Host "Done applying edits."
